$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '37.773.64'
$ws.Range("E2").Value = '  +5.83%  '
$ws.Range("D3").Value = '2.050.05'
$ws.Range("E3").Value = '  +3.29%  '
$ws.Range("D5").Value = '252.81'
$ws.Range("E5").Value = '  +4.24%  '
$ws.Range("D6").Value = '0.652'
$ws.Range("E6").Value = '  +2.09%  '
$ws.Range("D7").Value = '65.56'
$ws.Range("E7").Value = '  +15.04%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +5.03%  '
$ws.Range("D10").Value = '59.88'
$ws.Range("E10").Value = '  -0.23%  '
$ws.Range("D11").Value = '0.0762'
$ws.Range("E11").Value = '  +4.38%  '
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").Value = '0.926'
$ws.Range("E13").Value = '  +0.13%  '
$ws.Range("D14").Value = '15.23'
$ws.Range("E14").Value = '  +8.00%  '
$ws.Range("D15").Value = '2.351.30'
$ws.Range("E15").Value = '  +3.35%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '20.74'
$ws.Range("E16").Value = '  +20.24%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = '5.57'
$ws.Range("E17").Value = '  +5.95%  '
$ws.Range("D18").Value = '2.035.38'
$ws.Range("E18").Value = '  +2.55%  '
$ws.Range("D19").Value = '37.674.00'
$ws.Range("E19").Value = '  +6.03%  '
$ws.Range("D20").Value = '74.02'
$ws.Range("E20").Value = '  +5.06%  '
$ws.Range("D21").Value = '0.0₃0879'
$ws.Range("E21").Value = '  +4.90%  '
$ws.Range("D22").Value = '5.36'
$ws.Range("E22").Value = '  +5.82%  '
$ws.Range("D23").Value = '238.69'
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("E24").Value = '  +14.05%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  +3.99%  '
$ws.Range("D27").Value = '9.63'
$ws.Range("E27").Value = '  +5.51%  '
$ws.Range("D28").Value = '160.64'
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").Value = '19.97'
$ws.Range("E29").Value = '  +2.38%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '0.122'
$ws.Range("E30").Value = '  +2.48%  '
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '0.114'
$ws.Range("E31").Value = '  +27.16%  '
$ws.Range("E32").Value = '  +8.70%  '
$ws.Range("E33").Value = '  +6.71%  '
$ws.Range("D34").Value = '4.75'
$ws.Range("E34").Value = '  +11.55%  '
$ws.Range("D35").Value = '0.0619'
$ws.Range("E35").Value = '  +5.28%  '
$ws.Range("D36").Value = '2.43'
$ws.Range("E36").Value = '  +3.45%  '
$ws.Range("D37").Value = '1.85'
$ws.Range("E37").Value = '  +3.14%  '
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").Value = '6.11'
$ws.Range("E39").Value = '  +23.79%  '
$ws.Range("D40").Value = '0.103'
$ws.Range("E40").Value = '  +16.04%  '
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  +26.28%  '
$ws.Range("E42").Value = '  +4.21%  '
$ws.Range("D43").Value = '0.0219'
$ws.Range("E43").Value = '  +4.56%  '
$ws.Range("D44").Value = '2.92'
$ws.Range("E45").Value = '  +5.39%  '
$ws.Range("D46").Value = '17.07'
$ws.Range("E46").Value = '  +10.58%  '
$ws.Range("D47").Value = '8.01'
$ws.Range("E47").Value = '  +7.89%  '
$ws.Range("D48").Value = '95.28'
$ws.Range("E48").Value = '  +4.74%  '
$ws.Range("D49").Value = '1.416.14'
$ws.Range("E49").Value = '  +2.70%  '
$ws.Range("D50").Value = '2.94'
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("D51").Value = '47.47'
$ws.Range("E51").Value = '  +3.82%  '
